$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the Fecha/Volumen/Precio values between row 2 <-> row 8, and row 4 <-> row 5.
# Columns involved: D (Fecha), M (Volumen), N (Precio minimo), O (Precio maximo),
# P (Precio promedio ponderado), S (Precio $/Kg)

$cols = @("D", "M", "N", "O", "P", "S")

function Swap-Rows($ws, $rowA, $rowB, $cols) {
    foreach ($col in $cols) {
        $rangeA = $ws.Range("$col$rowA")
        $rangeB = $ws.Range("$col$rowB")
        $valA = $rangeA.Value2
        $valB = $rangeB.Value2
        $rangeA.Value = $valB
        $rangeB.Value = $valA
    }
}

Swap-Rows $ws 2 8 $cols
Swap-Rows $ws 4 5 $cols
